$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# New pagination header columns
$ws.Range("C1").Value = "languages pagination"
$ws.Range("D1").Value = "tags pagination"
$ws.Range("E1").Value = "notes pagination"
$ws.Range("F1").Value = "vocabulary pagination"

# New pagination values
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 2

# Make "settings" the active sheet / tab, with F2 selected
$ws.Activate() | Out-Null
$ws.Range("F2").Select() | Out-Null
